# "fixed up pins and drop down spec menu"
# Rebuild the library_specialization (column I) values on Sheet1: rows 3-27 get
# assigned a distinct specialization tag (the drop-down menu options), rows
# that previously just inherited the "Little Library"/"Public Library"
# placeholder and had no real specialization are cleared, and the public
# library rows (34-52) are all tagged "Official Community".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value  = "Animals or Pets"
$ws.Range("I4").Value  = "Art"
$ws.Range("I5").Value  = "Automotive"
$ws.Range("I6").Value  = "Baby"
$ws.Range("I7").Value  = "Children"
$ws.Range("I8").Value  = "Cooking"
$ws.Range("I9").Value  = "Fantasy"
$ws.Range("I10").Value = "Fiction"
$ws.Range("I11").Value = "Finance"
$ws.Range("I12").Value = "Fitness"
$ws.Range("I13").Value = "History"
$ws.Range("I14").Value = "Home Improvement"
$ws.Range("I15").Value = "Horror"

# These two rows no longer carry a specialization tag at all (cell + its
# style are removed entirely, matching the column's default formatting).
$ws.Range("I16").Clear()
$ws.Range("I17").Clear()

$ws.Range("I18").Value = "Nature"
$ws.Range("I19").Value = "Religious"
$ws.Range("I20").Value = "Romance"
$ws.Range("I21").Value = "Satire"
$ws.Range("I22").Value = "Science Fiction"
$ws.Range("I23").Value = "Self Improvement"
$ws.Range("I24").Value = "Space"
$ws.Range("I25").Value = "Sports"
# NB: I27 is set before I26 so new shared strings are interned in the same
# order the original edit produced (Young Adults before Technical or
# Textbooks) even though row 26 precedes row 27 on the sheet.
$ws.Range("I27").Value = "Young Adults"
$ws.Range("I26").Value = "Technical or Textbooks"

$ws.Range("I28").Clear()

# This row keeps its (distinct) cell formatting but loses its value.
$ws.Range("I29").ClearContents()

$ws.Range("I30").Clear()
$ws.Range("I31").Clear()
$ws.Range("I32").Clear()
$ws.Range("I33").Clear()

# Public-library rows all become "Official Community".
$ws.Range("I34:I52").Value = "Official Community"

# Window / view tidy-up (scroll position + zoom) from the spec clean-up pass.
$win = $ws.Application.ActiveWindow
$win.Zoom = 120
$ws.Range("H51").Select()
